$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Author cell: drop the second co-author's name.
$authorCell = $t.Rows.Item(4).Cells.Item(2)
$authorCell.Range.Text = "Robertino Bristiel"

# 2) "Curso Basico" step table: the row between step "6" and the old step "7"
#    is missing its step number. Fill it in with "7" using the same
#    formatting (10pt / 10pt complex-script) as the other step numbers.
$numberCell = $t.Rows.Item(18).Cells.Item(1)
$numberCell.Range.Text = "7"
$numberCell.Range.Font.Size = 10
$numberCell.Range.Font.SizeBi = 10

# 3) The following row's step number ("7") must become "8" to keep the
#    sequence consistent with the newly inserted step 7.
$followingNumberCell = $t.Rows.Item(19).Cells.Item(1)
$followingNumberCell.Range.Text = "8"
